$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sibirev I. V.")

# Update the four grade cells in row 31 from 4 to 5
$ws.Range("C31").Value = 5
$ws.Range("D31").Value = 5
$ws.Range("E31").Value = 5
$ws.Range("F31").Value = 5

# Recalculate so the shared-formula total in K31 picks up the new values
$excel.Calculate()

# Update the active selection to match the author's final cursor position
$ws.Range("G31").Select()
